# Reorder the "Requisitos" entries so that the LOM3229 line moves from the
# first position (row 23) to the last position (row 25), shifting the
# other two entries (LOB1021, LOM3016) up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lom3229 = "LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)`n"
$lob1021 = "LOB1021 -  Física IV  (Requisito)`n"
$lom3016 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

$ws.Range("B23").Value = $lob1021
$ws.Range("C23").Value = $lob1021

$ws.Range("B24").Value = $lom3016
$ws.Range("C24").Value = $lom3016

$ws.Range("B25").Value = $lom3229
$ws.Range("C25").Value = $lom3229
